# Applies the "Doing Updates for Financials" commit:
#  - Many historical financial figures on the STE sheet are replaced with "NA"
#    placeholders (years for which data isn't available any more).
#  - One figure (Other Liabilities, column G, row 62) is corrected to a new
#    numeric value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where the full D:J data range becomes "NA"
$naRows = @(8, 9, 10, 12, 14, 15, 17, 18, 20, 21, 22, 23, 24, 26, 27, 32, 33, 35, 81)

foreach ($r in $naRows) {
    $ws.Range("D" + $r + ":J" + $r).Value = "NA"
}

# Row 29 only had its D column populated with a number previously; it also
# becomes "NA" (columns E:J on that row were already "NA").
$ws.Range("D29").Value = "NA"

# Other Liabilities (row 62), column G: value correction 126300 -> 119400
$ws.Range("G62").Value = 119400
